$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.159.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.63%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.170.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -8.00%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'561.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.90%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'170.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.03%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.68%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.166.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -8.04%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -5.74%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.52%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.394"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.721.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -7.95%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.97%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'27.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.24%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.132.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.66%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -5.07%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.171.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -7.63%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.10%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'351.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.54%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'7.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -5.71%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.09%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'68.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.52%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.502"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.23%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -3.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.57%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -2.07%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.58%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.12%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.03%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.35%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'22.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.75%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.88%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -6.84%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'155.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.17%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.25%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.805"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -8.56%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'25.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -10.26%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Stacks"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'1.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.99%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'dogwifhat"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'2.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.650.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.84%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'RenderToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'6.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.47%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Filecoin"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'4.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -6.76%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Bittensor"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'328.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.91%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Hedera"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.0650"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.89%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'38.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.25%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -4.10%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0270"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -6.43%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.44%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.02%  "
$ws.Range("E51").Style = "Normal"
